# This edit removes the "Palestine" row from the origin/destination country
# list on the worksheet. Deleting the entire row shifts every row below it
# up by one (carrying along each cell's value/style), which is exactly the
# change captured by the diff: the shared string "Palestine" disappears and
# all the countries that were below it move up one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row that contains "Palestine" in column A and delete it entirely.
$palestineCell = $ws.Cells.Item(43, 1)
if ($palestineCell.Text -eq "Palestine") {
    $palestineRow = $palestineCell.Row
} else {
    $palestineRow = $ws.Range("A1:A200").Find("Palestine").Row
}

$ws.Rows.Item($palestineRow).Delete()

# Restore the view/selection state recorded in the saved workbook: the user
# had scrolled so row 7 is at the top and selected the entire row that now
# holds "Romania" (row 43 after the deletion).
$win = $excel.ActiveWindow
$excel.Goto($ws.Range("A7"), $true)
$ws.Range("A43:XFD43").Select()
